# Add a new "p_reported" column (reported p-values) to the study-results sheet.
# This inserts a new column at J (shifting the existing J/K/L "N" / "published" /
# "Notes" columns one place to the right, to K/L/M), fills in the header and the
# per-row reported-p-value labels, wordsmithing the dataset per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before J; N/published/Notes (J:L) shift right to K:M.
$ws.Columns("J:J").Insert()

# New column header.
$ws.Range("J1").Value = "p_reported"

# Populate p_reported per row. Values are written in this specific order so the
# shared-string table is built up in the same sequence as the source edit
# (first occurrence of each distinct label: >0.1, <0.1, <0.01, <0.05).
$ws.Range("J4").Value = ">0.1"
$ws.Range("J2").Value = "<0.1"
$ws.Range("J3").Value = "<0.01"
$ws.Range("J8").Value = "<0.05"

$ws.Range("J5").Value = ">0.1"
$ws.Range("J6").Value = ">0.1"
$ws.Range("J7").Value = ">0.1"
$ws.Range("J9").Value = "<0.05"
$ws.Range("J10").Value = "<0.05"
$ws.Range("J11").Value = ">0.1"
$ws.Range("J12").Value = ">0.1"
$ws.Range("J13").Value = "<0.05"
$ws.Range("J14").Value = "<0.01"
$ws.Range("J15").Value = "<0.01"
$ws.Range("J16").Value = "<0.01"
$ws.Range("J17").Value = "<0.01"
$ws.Range("J18").Value = "<0.01"
$ws.Range("J19").Value = "<0.01"
$ws.Range("J20").Value = "<0.05"
$ws.Range("J21").Value = "<0.01"
